$d = $word.ActiveDocument

# 1. Fix typo "Des encriptación" -> "Desencriptación"
$d.Content.Find.Execute("Des encriptación", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Desencriptación", 2)

# 2-4. The paragraphs "Se enviarán datos de la forma N" (N = 2,3,4) currently
# have the leading text "Se enviarán datos de la forma " and the digit "N"
# split across two separate runs, followed by a third run " (pN)". Merge only
# the first two runs into one, leaving the " (pN)" run untouched.
#
# This engine consolidates every run in a paragraph that shares identical
# formatting whenever any text inside the paragraph is edited. To keep the
# trailing " (pN)" run separate we briefly give it distinguishing
# (Bold) formatting before editing, then clear it again afterwards -
# toggling Bold off fully removes the property again, leaving no residue.
foreach ($n in 2,3,4) {
    $final = "Se enviarán datos de la forma $n"
    $marker = " (p$n)"

    # Locate " (pN)" and temporarily bold it so it keeps its own run.
    $mRange = $d.Content
    $mRange.Find.Execute($marker, $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
    $protect = $d.Range($mRange.Start, $mRange.End)
    $protect.Font.Bold = $true

    # Force the "forma " + "N" runs to merge by routing the text through a
    # placeholder (a same-text assignment is treated as a no-op and would
    # leave the runs split).
    $target = $d.Content
    $target.Find.Execute($final, $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
    $rng = $d.Range($target.Start, $target.End)
    $rng.Text = "__TMP__"
    $d.Content.Find.Execute("__TMP__", $true, $false, $false, $false, $false,
                             $true, 1, $false, $final, 2)

    # Un-bold the " (pN)" run again (fully clears the property).
    $restore = $d.Content
    $restore.Find.Execute($marker, $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
    $unprotect = $d.Range($restore.Start, $restore.End)
    $unprotect.Font.Bold = $false
}
